# Course Management UI update:
#  - add "등록 강좌" (registered courses) and "비밀번호" (password) columns
#  - add a first student record (염승욱) with course list + credit count + password
#  - bold the "학점" / "등록 강좌" headers and the student-name cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "학번"
$ws.Range("B1").Value = "이름"
$ws.Range("C1").Value = "학과"
$ws.Range("D1").Value = "주민등록번호"
$ws.Range("E1").Value = "학점"
$ws.Range("F1").Value = "등록 강좌"
$ws.Range("G1").Value = "비밀번호"

# ---- Data row (row 2) ----
$ws.Range("A2").Value = 123
$ws.Range("B2").Value = "염승욱"
$ws.Range("D2").Value = 123456789
$ws.Range("E2").Value = 18
$ws.Range("F2").Value = "자바,자바2,자바3"
$ws.Range("G2").Value = "thdjs07"

# ---- Emphasis on the course-load related headers / new name ----
$ws.Range("E1").Font.Bold = $true
$ws.Range("F1").Font.Bold = $true
$ws.Range("B2").Font.Bold = $true

# ---- Selection moves to the new "password" header cell ----
$ws.Range("G1").Select() | Out-Null
